$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45231 to 45232
# for every data row (rows 2 through 536).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 536) { $lastRow = 536 }

$ws.Range("C2:C$lastRow").Value = 45232
